$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.007.33"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "1.553.61"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'286.41"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.3773"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").Value = "'0.3236"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "'1.121"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").Value = "'41.20"
$ws.Range("E10").Value = "  -12.82%  "
$ws.Range("D11").Value = "'0.07298"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "'19.34"
$ws.Range("E13").Value = "  -6.80%  "
$ws.Range("D14").Value = "'5.710"
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("D15").Value = "'6.799"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "1.554.31"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "'0.00001079"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "'0.06649"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'85.03"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("D20").Value = "'6.426"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'15.96"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").Value = "'11.51"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").Value = "22.040.36"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "'2.271"
$ws.Range("E25").Value = "  -4.54%  "
$ws.Range("D26").Value = "'2.510"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("D27").Value = "'149.97"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "'18.83"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").Value = "'4.843"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("D30").Value = "1.725.24"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "'119.99"
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").Value = "'1.121"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("D33").Value = "'5.906"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("D34").Value = "'0.08164"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").Value = "'9.277"
$ws.Range("E35").Value = "  -6.37%  "
$ws.Range("D36").Value = "'1.642"
$ws.Range("E36").Value = "  -17.49%  "
$ws.Range("D37").Value = "'5.226"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("E38").Value = "  -6.54%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").Value = "'0.2109"
$ws.Range("E40").Value = "  -4.87%  "
$ws.Range("D41").Value = "'1.215"
$ws.Range("E41").Value = "  -6.77%  "
$ws.Range("D42").Value = "'10.87"
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'0.5939"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("D45").Value = "'13.55"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").Value = "'3.724"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").Value = "'0.5727"
$ws.Range("E47").Value = "  -5.29%  "
$ws.Range("D48").Value = "'1.934"
$ws.Range("E48").Value = "  -5.21%  "
$ws.Range("D49").Value = "'119.80"
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("D50").Value = "'1.153"
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("D51").Value = "'0.06902"
$ws.Range("E51").Value = "  -4.07%  "
